$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-11-17 Monday" "2025-11-18 Tuesday"

Replace-Text "788×4=3152" "842×4=3368"
Replace-Text "340×6=2040" "290×7=2030"
Replace-Text "746×4=2984" "624×6=3744"
Replace-Text "120×3=360" "744×9=6696"
Replace-Text "785×2=1570" "874×8=6992"

Replace-Text "392×2=784" "501×7=3507"
Replace-Text "180×9=1620" "232×8=1856"
Replace-Text "324×2=648" "403×8=3224"
Replace-Text "257×6=1542" "365×3=1095"
Replace-Text "268×4=1072" "802×6=4812"

Replace-Text "895×9=8055" "871×9=7839"
Replace-Text "713×9=6417" "737×8=5896"
Replace-Text "910×5=4550" "750×7=5250"
Replace-Text "646×7=4522" "411×4=1644"
Replace-Text "325×5=1625" "551×6=3306"

Replace-Text "400×4=1600" "896×3=2688"
Replace-Text "903×8=7224" "253×7=1771"
Replace-Text "210×9=1890" "637×8=5096"
Replace-Text "345×7=2415" "255×7=1785"
Replace-Text "971×8=7768" "207×7=1449"

Replace-Text "339×8=2712" "282×6=1692"
Replace-Text "939×4=3756" "885×6=5310"
Replace-Text "553×2=1106" "338×4=1352"
Replace-Text "831×6=4986" "470×2=940"
Replace-Text "312×4=1248" "104×2=208"
